# Hortaliza, Vega Modelo de Temuco - Ciboulette
# A new weekly price record is inserted as row 60 (pushing the existing
# rows 60-151 down to 61-152, which is exactly what the diff shows: each
# old row's data reappears one row lower, and a brand-new row 152 carries
# what used to be the last row, 151).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 60; everything below shifts down one.
$ws.Rows(60).Insert()

# Populate the newly inserted row 60 with the new observation.
$ws.Range("A60").Value = 10
$ws.Range("B60").Value = "Vega Modelo de Temuco"
$ws.Range("C60").Value = "La Araucanía"
$ws.Range("D60").Value = 44467
$ws.Range("E60").Value = 9
$ws.Range("F60").Value = 100112039
$ws.Range("G60").Value = "Ciboulette"
$ws.Range("H60").Value = "Sin especificar"
$ws.Range("I60").Value = "Primera"
$ws.Range("J60").Value = 20
$ws.Range("K60").Value = 7000
$ws.Range("L60").Value = 7000
$ws.Range("M60").Value = 7000
$ws.Range("N60").Value = "$/docena de atados"
$ws.Range("O60").Value = "Provincia de Cautín"
$ws.Range("P60").Value = 2333
$ws.Range("Q60").Value = 3
$ws.Range("R60").Value = "Hortaliza"
